# Updates cryptos list values (Price + Volume(1h) columns) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.281.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.619.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.487"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0615"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0816"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.844.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.618.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.286.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("E30").Value = "  +10.52%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.181.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.16%  "
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.805"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.10%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.496"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.790"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("E43").Value = "  +4.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.755.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("E46").Value = "  +14.22%  "
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.409"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("E51").Value = "  -0.30%  "
